# PR Stamina - Combat (#121)
# Adds a new worksheet "Arkusz6" with stamina/hp stat tables, wired up with
# formulas derived from class base stats (str/end/dex) previously recorded
# on Arkusz5.

$wb = $excel.ActiveWorkbook

# --- add the new sheet as the last tab (becomes the active tab, matching
#     the workbook's activeTab bump from 4 -> 5) -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Arkusz6"

# --- shared-string write order mirrors the target file so new <si> entries
#     land at indices 24..34 in the same order as the authored workbook -----
$ws.Range("A1").Value = "level"
$ws.Range("A3").Value = "end"
$ws.Range("A9").Value = "base str"
$ws.Range("A10").Value = "base end"
$ws.Range("A11").Value = "base dex"
$ws.Range("C10").Value = "bonus end"
$ws.Range("C9").Value = "bonus str"
$ws.Range("C11").Value = "bonus dex"
$ws.Range("A6").Value = "hp"
$ws.Range("A7").Value = "stamina"
$ws.Range("A14").Value = "hp v"
$ws.Range("A2").Value = "str"
$ws.Range("A4").Value = "dex"

# --- row 1: level header values ---------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 10
$ws.Range("E1").Value = 15
$ws.Range("F1").Value = 20

# --- row 9-11: base class attributes + bonus per 5 levels -------------------
$ws.Range("B9").Value = 65
$ws.Range("D9").Value = 6.25

$ws.Range("B10").Value = 65
$ws.Range("D10").Value = 6.25

$ws.Range("B11").Value = 55
$ws.Range("D11").Value = 3.75

# --- row 2: str per level ----------------------------------------------------
$ws.Range("B2").Formula = '=$B9+$D9*B$1/5'
$ws.Range("C2:F2").Formula = '=$B9+$D9*C$1/5'

# --- row 3: end per level ----------------------------------------------------
$ws.Range("B3:F3").Formula = '=$B10+$D10*B$1/5'

# --- row 4: dex per level ----------------------------------------------------
$ws.Range("B4:F4").Formula = '=$B11+$D11*B$1/5'

# --- row 6: hp per level ------------------------------------------------------
$ws.Range("B6").Formula = '=500 * (1 + (B14-50)/50)'
$ws.Range("C6:E6").Formula = '=500 * (1 + (C14-50)/50)'
$ws.Range("F6").Formula = '=500 * (1 + (F14-50)/50)'

# --- row 7: stamina per level --------------------------------------------------
$ws.Range("B7").Formula = '=50+B3*2.5+B4*2'
$ws.Range("C7:F7").Formula = '=50+C3*2.5+C4*2'

# --- row 14: hp v (weighted avg of end/str) used by row 6 -------------------
$ws.Range("B14").Formula = '=B3*0.8 + B2*0.2'
$ws.Range("C14:F14").Formula = '=C3*0.8 + C2*0.2'

# --- sheet view: selection + it becomes the active/visible tab --------------
$ws.Range("I10").Select()

# --- Arkusz4's remembered selection also moved (B11) in the authored edit ---
$ws4 = $wb.Worksheets.Item("Arkusz4")
$ws4.Range("B11").Select()

# --- re-select the new sheet so it ends up the active one after touching
#     Arkusz4's selection above ----------------------------------------------
$ws.Select()
$ws.Range("I10").Select()
